$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reorder the comma-separated "Recorded By" values in column G:
# move the first name/email in the list to the end (left-rotate by one).
$rows = @(2, 3, 4, 5, 6, 11, 12, 13, 14, 15, 29, 30, 32, 33, 38, 39, 40, 41, 42, 56, 57, 58, 59, 60, 65, 66, 67, 68, 69, 84, 85, 86, 89, 90, 93, 95, 110, 111, 112, 115, 116, 119, 121, 136, 137, 138, 141, 142, 145, 147)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $parts = $cell.Value2 -split ", "
    if ($parts.Count -gt 1) {
        $newParts = $parts[1..($parts.Count - 1)] + $parts[0]
        $cell.Value2 = [string]::Join(", ", $newParts)
    }
}
